# Apply "owner pdf link send to the client" update:
# Update OFF ROAD DAY (L) and Rate (O) values for several vendor rows,
# restore row 22's Rate/Amount back to the standard values, and move the
# active selection/scroll position in the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Updated OFF ROAD DAY (L) and Rate (O) figures ---
$ws.Range("L2").Value = 1
$ws.Range("O2").Value = 600

$ws.Range("L4").Value = 2
$ws.Range("O4").Value = 643.23

$ws.Range("L5").Value = 4
$ws.Range("O5").Value = 546.54

$ws.Range("O6").Value = 246.45

$ws.Range("O7").Value = 600

$ws.Range("O8").Value = 450

$ws.Range("L9").Value = 2
$ws.Range("O9").Value = 746.45

$ws.Range("L10").Value = 4

$ws.Range("L11").Value = 5
$ws.Range("O11").Value = 494.56

$ws.Range("L12").Value = 1
$ws.Range("O12").Value = 556.66999999999996

$ws.Range("L13").Value = 2

$ws.Range("L17").Value = 2
$ws.Range("O17").Value = 500

$ws.Range("L18").Value = 4
$ws.Range("O18").Value = 600

$ws.Range("O19").Value = 756.65

$ws.Range("O20").Value = 459.69

$ws.Range("L21").Value = 5

# Row 22 reverted back to the standard rate/amount
$ws.Range("O22").Value = 693.54
$ws.Range("P22").Value = 21500

$ws.Range("L24").Value = 6
$ws.Range("O24").Value = 600

$ws.Range("O25").Value = 456.56

$ws.Range("L26").Value = 3
$ws.Range("O26").Value = 345.56

$ws.Range("L28").Value = 1

$ws.Range("O29").Value = 456.59

$ws.Range("L31").Value = 1

$ws.Range("O33").Value = 334.56

$ws.Range("L35").Value = 1

$ws.Range("L38").Value = 2

$ws.Range("L40").Value = 3

# --- Scroll / selection state (topLeftCell -> J16, selection -> O33) ---
$win = $wb.Windows.Item(1)
$win.ScrollRow = 16
$win.ScrollColumn = 10
$ws.Range("O33").Select()
